$wb = $excel.ActiveWorkbook

$wsNav = $wb.Worksheets.Item("NAV")
$wsNav.Range("A6").Value = "ga4_path"
$wsNav.Range("B6").Value = "/en-us/nav-ptv/nav/4-passenger/nav-4e"
$wsNav.Activate()
$wsNav.Range("B6").Select()

$ws5525 = $wb.Worksheets.Item("5525")
$ws5525.Range("A6").Value = "ga4_path"
$ws5525.Range("B6").Value = "en-us/side-x-side/teryx/teryx-4-5-passenger-supercharged/teryx4-5-h2"
$ws5525.Activate()
$ws5525.Range("B6").Select()
